$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name text correction (B5) ---
$ws.Range("B5").Value = "کیمیا-ص. معدنی کیمیای زنجان گستران"

# --- Header row 8: financial period labels (shift left, append newest period) ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates (shift left, append newest date) ---
$ws.Range("D9").Value = "1400-11-11 (4)"
$ws.Range("E9").Value = "1401-04-08 (9)"
$ws.Range("F9").Value = "1401-05-11 (4)"
$ws.Range("G9").Value = "1401-08-29 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-28 (8)"
$ws.Range("J9").Value = "1401-05-11 (2)"
$ws.Range("K9").Value = "1401-08-29 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-28"

# --- Data rows 11-27: shift each quarterly column left, append newest period values ---
# Row 11
$ws.Range("D11").Value = 1579946
$ws.Range("E11").Value = 4246792
$ws.Range("F11").Value = 1095405
$ws.Range("G11").Value = 3259628
$ws.Range("H11").Value = 4308768
$ws.Range("I11").Value = 7078947
$ws.Range("J11").Value = 1416204
$ws.Range("K11").Value = 4139043
$ws.Range("L11").Value = 5294306
$ws.Range("M11").Value = 8846679

# Row 12
$ws.Range("D12").Value = -976752
$ws.Range("E12").Value = -2759231
$ws.Range("F12").Value = -711157
$ws.Range("G12").Value = -2117016
$ws.Range("H12").Value = -2839940
$ws.Range("I12").Value = -4991386
$ws.Range("J12").Value = -982683
$ws.Range("K12").Value = -2780908
$ws.Range("L12").Value = -3532249
$ws.Range("M12").Value = -5785509

# Row 13
$ws.Range("D13").Value = 603194
$ws.Range("E13").Value = 1487561
$ws.Range("F13").Value = 384248
$ws.Range("G13").Value = 1142612
$ws.Range("H13").Value = 1468828
$ws.Range("I13").Value = 2087561
$ws.Range("J13").Value = 433521
$ws.Range("K13").Value = 1358135
$ws.Range("L13").Value = 1762057
$ws.Range("M13").Value = 3061170

# Row 14
$ws.Range("D14").Value = -34860
$ws.Range("E14").Value = -73858
$ws.Range("F14").Value = -24333
$ws.Range("G14").Value = -134648
$ws.Range("H14").Value = -183933
$ws.Range("I14").Value = -317017
$ws.Range("J14").Value = -50408
$ws.Range("K14").Value = -161360
$ws.Range("L14").Value = -242065
$ws.Range("M14").Value = -362437

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 10500
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 10500
$ws.Range("H16").Value = 15750
$ws.Range("I16").Value = 21000
$ws.Range("J16").Value = 7350
$ws.Range("K16").Value = 14700
$ws.Range("L16").Value = 22050
$ws.Range("M16").Value = 923046

# Row 17
$ws.Range("D17").Value = 568334
$ws.Range("E17").Value = 1424203
$ws.Range("F17").Value = 359915
$ws.Range("G17").Value = 1018464
$ws.Range("H17").Value = 1300645
$ws.Range("I17").Value = 1791544
$ws.Range("J17").Value = 390463
$ws.Range("K17").Value = 1211475
$ws.Range("L17").Value = 1542042
$ws.Range("M17").Value = 3621779

# Row 18
$ws.Range("D18").Value = -12430
$ws.Range("E18").Value = -14844
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = -3646
$ws.Range("J18").Value = -6045
$ws.Range("K18").Value = -10556
$ws.Range("L18").Value = -10556
$ws.Range("M18").Value = -12575

# Row 19
$ws.Range("D19").Value = 148188
$ws.Range("E19").Value = 148793
$ws.Range("F19").Value = 129285
$ws.Range("G19").Value = 276232
$ws.Range("H19").Value = 282355
$ws.Range("I19").Value = 282249
$ws.Range("J19").Value = 55880
$ws.Range("K19").Value = 159379
$ws.Range("L19").Value = 173958
$ws.Range("M19").Value = 186205

# Row 20
$ws.Range("D20").Value = 704092
$ws.Range("E20").Value = 1558152
$ws.Range("F20").Value = 489200
$ws.Range("G20").Value = 1294696
$ws.Range("H20").Value = 1583000
$ws.Range("I20").Value = 2070147
$ws.Range("J20").Value = 440298
$ws.Range("K20").Value = 1360298
$ws.Range("L20").Value = 1705444
$ws.Range("M20").Value = 3795409

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("D22").Value = 704092
$ws.Range("E22").Value = 1558152
$ws.Range("F22").Value = 489200
$ws.Range("G22").Value = 1294696
$ws.Range("H22").Value = 1583000
$ws.Range("I22").Value = 2070147
$ws.Range("J22").Value = 440298
$ws.Range("K22").Value = 1360298
$ws.Range("L22").Value = 1705444
$ws.Range("M22").Value = 3795409

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 704092
$ws.Range("E24").Value = 1558152
$ws.Range("F24").Value = 489200
$ws.Range("G24").Value = 1294696
$ws.Range("H24").Value = 1583000
$ws.Range("I24").Value = 2070147
$ws.Range("J24").Value = 440298
$ws.Range("K24").Value = 1360298
$ws.Range("L24").Value = 1705444
$ws.Range("M24").Value = 3795409

# Row 25
$ws.Range("D25").Value = 671
$ws.Range("E25").Value = 1484
$ws.Range("F25").Value = 245
$ws.Range("G25").Value = 647
$ws.Range("H25").Value = 792
$ws.Range("I25").Value = 1972
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 680
$ws.Range("L25").Value = 853
$ws.Range("M25").Value = 1265

# Row 26
$ws.Range("D26").Value = 1050000
$ws.Range("E26").Value = 1050000
$ws.Range("F26").Value = 2000000
$ws.Range("G26").Value = 2000000
$ws.Range("H26").Value = 2000000
$ws.Range("I26").Value = 1050000
$ws.Range("J26").Value = 2000000
$ws.Range("K26").Value = 2000000
$ws.Range("L26").Value = 2000000
$ws.Range("M26").Value = 3000000

# Row 27
$ws.Range("D27").Value = 235
$ws.Range("E27").Value = 519
$ws.Range("F27").Value = 163
$ws.Range("G27").Value = 432
$ws.Range("H27").Value = 528
$ws.Range("I27").Value = 690
$ws.Range("J27").Value = 147
$ws.Range("K27").Value = 453
$ws.Range("L27").Value = 568
$ws.Range("M27").Value = 1265
